# CI: Auto Update Data — refresh MAA operator proficiency stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "maa://24702 (94.1), maa://25390 (97.01), maa://36681 (90.77)"
$ws.Range("AA2").Value = "maa://21246 (91.2), maa://36684 (98.63), ***maa://22731 (6.67)"
$ws.Range("AE2").Value = "maa://25251 (92.5), ***maa://21730 (17.19), ***maa://39501 (21.43), *maa://36675 (60.0)"
$ws.Range("AA3").Value = "maa://24390 (96.08)"
$ws.Range("W4").Value = "**maa://32495 (47.54), ***maa://31785 (18.75), ***maa://36683 (26.67)"
$ws.Range("AA4").Value = "*maa://32658 (68.75)"
$ws.Range("C6").Value = "maa://42407 (87.5)"
$ws.Range("W7").Value = "maa://22399 (94.66), *maa://22758 (71.7)"
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "3"
$scratch.Copy()
$ws.Range("AD7").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("AE7").Value = "*maa://26191 (68.49), *maa://36671 (72.09), maa://42530 (100.0)"
$ws.Range("W8").Value = "maa://21411 (96.03)"
$ws.Range("AE10").Value = "*maa://25021 (56.16), *maa://22733 (58.62), maa://22761 (100.0)"
$ws.Range("W11").Value = "maa://36713 (97.85)"
$ws.Range("W12").Value = "maa://22753 (91.22), *maa://21485 (76.74), maa://37962 (81.25)"
$ws.Range("W13").Value = "*maa://34957 (76.6), *maa://22768 (53.33)"
$ws.Range("AE13").Value = "**maa://22737 (30.6), maa://39883 (88.46), *maa://39885 (73.68)"
$ws.Range("K14").Value = "maa://26245 (96.12), maa://21288 (96.21), maa://36682 (100.0), maa://39841 (93.18)"
$ws.Range("G15").Value = "maa://24304 (88.59), maa://21478 (91.18)"
$ws.Range("C16").Value = "maa://21441 (96.17), maa://36679 (94.12), maa://37650 (95.45)"
$ws.Range("W16").Value = "maa://28501 (97.47), maa://28051 (95.83)"
$ws.Range("C18").Value = "maa://24570 (96.61)"
$ws.Range("G18").Value = "maa://24421 (90.23)"
$ws.Range("S19").Value = "maa://24386 (98.81)"
$ws.Range("K20").Value = "maa://41331 (82.93)"
$ws.Range("O20").Value = "maa://37442 (96.55)"
$ws.Range("K23").Value = "maa://39756 (92.5), maa://39875 (95.83)"
$ws.Range("W24").Value = "maa://23504 (92.9), maa://29988 (86.27), **maa://22892 (40.14), *maa://25141 (77.05), maa://36663 (80.7), ***maa://22815 (23.08)"
$ws.Range("C25").Value = "maa://29753 (95.15)"
$ws.Range("AA25").Value = "maa://31215 (84.15), *maa://24516 (79.07), maa://26001 (88.89)"
$ws.Range("AA26").Value = "*maa://42235 (73.68)"
$ws.Range("S27").Value = "*maa://30624 (76.32)"
$ws.Range("C28").Value = "maa://24465 (90.36), maa://25725 (82.28)"
$ws.Range("W28").Value = "maa://39929 (86.7), ***maa://39723 (14.71), maa://41749 (81.25)"
$ws.Range("AE28").Value = "maa://36660 (93.85), *maa://36701 (64.0)"
$ws.Range("C29").Value = "maa://31694 (97.92)"
$ws.Range("K30").Value = "maa://30442 (94.55)"
$ws.Range("W30").Value = "*maa://39477 (75.0)"
$ws.Range("S32").Value = "maa://41108 (91.43), maa://41238 (94.59)"
$ws.Range("S34").Value = "maa://24526 (93.13)"
$ws.Range("K35").Value = "maa://41296 (98.11)"
$ws.Range("AE38").Value = "maa://36697 (84.62)"
$ws.Range("G39").Value = "maa://25199 (86.11), maa://36670 (88.06), maa://30434 (87.5), ***maa://25036 (16.0)"
$ws.Range("O39").Value = "maa://24709 (92.16)"
$ws.Range("G43").Value = "maa://22525 (92.68), maa://21284 (82.93)"
$ws.Range("G46").Value = "maa://35931 (92.58)"
$ws.Range("G55").Value = "maa://32532 (91.93)"

Write-Output "Applied CI data refresh to 46 cell(s)."
